$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "St6gal1"
$ws.Range("C2").Value = "Cd22"
$ws.Range("D2").Value = "M2"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.327185666666666
$ws.Range("H2").Value = 9.981556999999999
$ws.Range("I2").Value = 0.1584084717220711
$ws.Range("J2").Value = 0.1584084717220711
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 14.445741
$ws.Range("N2").Value = 43.337223
$ws.Range("O2").Value = 0.9879457199603994
$ws.Range("P2").Value = 0.9879457199603994
$ws.Range("Q2").Value = 48.06366239957899
$ws.Range("R2").Value = 432.572961596211
$ws.Range("S2").Value = 0.1564989716432881
$ws.Range("T2").Value = 0.1564989716432881

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "St6gal1"
$ws.Range("C3").Value = "Cd22"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.327185666666666
$ws.Range("H3").Value = 9.981556999999999
$ws.Range("I3").Value = 0.1584084717220711
$ws.Range("J3").Value = 0.1584084717220711
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.1762576666666667
$ws.Range("N3").Value = 0.528773
$ws.Range("O3").Value = 0.01205428003960061
$ws.Range("P3").Value = 0.01205428003960061
$ws.Range("Q3").Value = 0.5864419821734445
$ws.Range("R3").Value = 5.277977839561
$ws.Range("S3").Value = 0.001909500078782999
$ws.Range("T3").Value = 0.001909500078782998

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "St6gal1"
$ws.Range("C4").Value = "Cd22"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.442036
$ws.Range("H4").Value = 7.326108
$ws.Range("I4").Value = 0.1162661869236271
$ws.Range("J4").Value = 0.1162661869236271
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 14.445741
$ws.Range("N4").Value = 43.337223
$ws.Range("O4").Value = 0.9879457199603994
$ws.Range("P4").Value = 0.9879457199603994
$ws.Range("Q4").Value = 35.277019568676
$ws.Range("R4").Value = 317.493176118084
$ws.Range("S4").Value = 0.1148646817473132
$ws.Range("T4").Value = 0.1148646817473132

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "St6gal1"
$ws.Range("C5").Value = "Cd22"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.442036
$ws.Range("H5").Value = 7.326108
$ws.Range("I5").Value = 0.1162661869236271
$ws.Range("J5").Value = 0.1162661869236271
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1762576666666667
$ws.Range("N5").Value = 0.528773
$ws.Range("O5").Value = 0.01205428003960061
$ws.Range("P5").Value = 0.01205428003960061
$ws.Range("Q5").Value = 0.430427567276
$ws.Range("R5").Value = 3.873848105484
$ws.Range("S5").Value = 0.001401505176313951
$ws.Range("T5").Value = 0.001401505176313951

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "St6gal1"
$ws.Range("C6").Value = "Cd22"
$ws.Range("D6").Value = "M2"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 10.64260933333333
$ws.Range("H6").Value = 31.927828
$ws.Range("I6").Value = 0.5066983476510879
$ws.Range("J6").Value = 0.5066983476510878
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 14.445741
$ws.Range("N6").Value = 43.337223
$ws.Range("O6").Value = 0.9879457199603994
$ws.Range("P6").Value = 0.9879457199603994
$ws.Range("Q6").Value = 153.740377993516
$ws.Range("R6").Value = 1383.663401941644
$ws.Range("S6").Value = 0.5005904638728987
$ws.Range("T6").Value = 0.5005904638728986

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "St6gal1"
$ws.Range("C7").Value = "Cd22"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 10.64260933333333
$ws.Range("H7").Value = 31.927828
$ws.Range("I7").Value = 0.5066983476510879
$ws.Range("J7").Value = 0.5066983476510878
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.1762576666666667
$ws.Range("N7").Value = 0.528773
$ws.Range("O7").Value = 0.01205428003960061
$ws.Range("P7").Value = 0.01205428003960061
$ws.Range("Q7").Value = 1.875841488338222
$ws.Range("R7").Value = 16.882573395044
$ws.Range("S7").Value = 0.006107883778189118
$ws.Range("T7").Value = 0.006107883778189116

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "St6gal1"
$ws.Range("C8").Value = "Cd22"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.592005666666666
$ws.Range("H8").Value = 13.776017
$ws.Range("I8").Value = 0.2186269937032139
$ws.Range("J8").Value = 0.2186269937032139
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 14.445741
$ws.Range("N8").Value = 43.337223
$ws.Range("O8").Value = 0.9879457199603994
$ws.Range("P8").Value = 0.9879457199603994
$ws.Range("Q8").Value = 66.33492453119899
$ws.Range("R8").Value = 597.014320780791
$ws.Range("S8").Value = 0.2159916026968993
$ws.Range("T8").Value = 0.2159916026968993

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "St6gal1"
$ws.Range("C9").Value = "Cd22"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.592005666666666
$ws.Range("H9").Value = 13.776017
$ws.Range("I9").Value = 0.2186269937032139
$ws.Range("J9").Value = 0.2186269937032139
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.1762576666666667
$ws.Range("N9").Value = 0.528773
$ws.Range("O9").Value = 0.01205428003960061
$ws.Range("P9").Value = 0.01205428003960061
$ws.Range("Q9").Value = 0.8093762041267777
$ws.Range("R9").Value = 7.284385837141
$ws.Range("S9").Value = 0.002635391006314539
$ws.Range("T9").Value = 0.002635391006314539

Write-Output "Edit complete"